# "Generate Report for Handback" - refresh the localization-status report:
#   - Status moves from "Ready for handoff" to "Handed back: in sync with en-US"
#   - The de-de handback is no longer stale, so its "Error Detail" explanation
#     is cleared (same for zh-cn, which already reported no error)
#   - Handback timestamps for zh-cn / de-de are refreshed
#   - A few report columns are widened now that the Error Detail text is gone

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Status: "Ready for handoff" -> "Handed back: in sync with en-US" ---
$overview.Range("E2").Value = "Handed back: in sync with en-US"
$overview.Range("F2").Value = "Handed back: in sync with en-US"
$zhcn.Range("C2").Value     = "Handed back: in sync with en-US"
$dede.Range("C2").Value     = "Handed back: in sync with en-US"

# --- zh-cn: refresh handback datetime, clear error detail ---
$zhcn.Range("K2").Value = "2016-08-31 09:01:58"
$zhcn.Range("P2").Value = ""

# --- de-de: refresh handback datetime, clear error detail ---
$dede.Range("K2").Value = "2016-08-31 09:02:19"
$dede.Range("P2").Value = ""

# --- Column width refresh (report columns widened) ---
$overview.Columns.Item(5).ColumnWidth = 29.1
$overview.Columns.Item(6).ColumnWidth = 29.1

$zhcn.Columns.Item(3).ColumnWidth  = 29.1
$zhcn.Columns.Item(16).ColumnWidth = 12.83

$dede.Columns.Item(3).ColumnWidth  = 29.1
$dede.Columns.Item(16).ColumnWidth = 12.83
